# MassQL_Queries.xlsx edit
#
# Summary of the change (per commit message / diff):
#  - Removed unused query-file columns: "M+H", "M-H", "Untargeted F5method7 RT",
#    "Untargeted F5method7 Comments", "INTEGRATION_MIN", "INTEGRATION_MAX".
#  - The RTMIN / RTMAX columns, which used to be formulas referencing the
#    (now removed) "Untargeted F5method7 RT" column, are "calculated" once
#    (flattened to static values) before that column disappears.
#  - The trailing "threshold" column header is renamed to "QC_threshold".
#  - The conditional formatting (duplicate-value highlighting) that lived on
#    the removed "M+H" column is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Freeze every formula in the used range to its last computed value.
#    This "calculates" the M+H / M-H / RTMIN / RTMAX formula columns down to
#    plain numbers before the columns they depend on are deleted below.
$used = $ws.UsedRange
$used.Value = $used.Value()

# 2) Drop the conditional formatting that was applied to the M+H column
#    (duplicate-value highlighting) -- that column is going away.
$ws.Cells.FormatConditions.Delete()

# 3) Remove the unused columns.
#    Original layout: A Name | B KEGG | C Formula | D Monoisotopic |
#      E M+H | F M-H | G Untargeted F5method7 RT | H ion_mode |
#      I TOLERANCEPPM | J RTMIN | K RTMAX | L INTEGRATION_MIN |
#      M INTEGRATION_MAX | N Untargeted F5method7 Comments | O threshold
#    Delete E:G (M+H, M-H, Untargeted F5method7 RT) ...
$ws.Columns("E:G").Delete()
#    ... which shifts everything left by 3, so INTEGRATION_MIN/MAX and
#    Untargeted F5method7 Comments are now at I:K -- delete those too.
$ws.Columns("I:K").Delete()

# Final layout: A Name | B KEGG | C Formula | D Monoisotopic | E ion_mode |
#   F TOLERANCEPPM | G RTMIN | H RTMAX | I QC_threshold (renamed)

# 4) Rename the trailing "threshold" header to "QC_threshold".
$ws.Range("I1").Value = "QC_threshold"
